$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.957004961598832
$ws.Range("C2").Value = 4.313327511860538
$ws.Range("E2").Value = 18.02067970953283
$ws.Range("F2").Value = 44.02529512794491
$ws.Range("G2").Value = 42.20406576581988
$ws.Range("H2").Value = 17.57655709645851
$ws.Range("J2").Value = 9.003186998171328
$ws.Range("K2").Value = 9.348048149965148
$ws.Range("M2").Value = 17.1337054571114
$ws.Range("B3").Value = 9.745712650386167
$ws.Range("C3").Value = 4.16530387558455
$ws.Range("E3").Value = 17.95242954479495
$ws.Range("F3").Value = 43.91291568152698
$ws.Range("G3").Value = 42.14947246247466
$ws.Range("H3").Value = 17.61141076735169
$ws.Range("J3").Value = 9.021944783352584
$ws.Range("K3").Value = 9.219275531925582
$ws.Range("M3").Value = 17.0583100539414
$ws.Range("B4").Value = 9.616291574269566
$ws.Range("C4").Value = 4.070704107893452
$ws.Range("E4").Value = 17.91411173637775
$ws.Range("F4").Value = 43.85428624938123
$ws.Range("G4").Value = 42.12713406041926
$ws.Range("H4").Value = 17.63566683661126
$ws.Range("J4").Value = 9.034054532392327
$ws.Range("K4").Value = 9.141582228360626
$ws.Range("M4").Value = 17.01537319038351
$ws.Range("B5").Value = 9.563712625158391
$ws.Range("C5").Value = 4.031251148835703
$ws.Range("E5").Value = 17.8994106861418
$ws.Range("F5").Value = 43.83301616633935
$ws.Range("G5").Value = 42.12084278927706
$ws.Range("H5").Value = 17.64626823046671
$ws.Range("J5").Value = 9.039138713022188
$ws.Range("K5").Value = 9.110307226818568
$ws.Range("M5").Value = 16.9987338596876
$ws.Range("B6").Value = 9.554994106659693
$ws.Range("C6").Value = 4.024646458105429
$ws.Range("E6").Value = 17.89702509298001
$ws.Range("F6").Value = 43.82964298859061
$ws.Range("G6").Value = 42.11996790404653
$ws.Range("H6").Value = 17.64807184142141
$ws.Range("J6").Value = 9.039991969991847
$ws.Range("K6").Value = 9.105138493301913
$ws.Range("M6").Value = 16.99602310321704
$ws.Range("B7").Value = 9.615581713412482
$ws.Range("C7").Value = 4.070175643555292
$ws.Range("E7").Value = 17.91390975864263
$ws.Range("F7").Value = 43.85398876164821
$ws.Range("G7").Value = 42.1270378312704
$ws.Range("H7").Value = 17.63580690998273
$ws.Range("J7").Value = 9.034122494154831
$ws.Range("K7").Value = 9.141158828545176
$ws.Range("M7").Value = 17.01514529592004
$ws.Range("B8").Value = 9.884135240329071
$ws.Range("C8").Value = 4.263076919854282
$ws.Range("E8").Value = 17.99640957474106
$ws.Range("F8").Value = 43.984403060971
$ws.Range("G8").Value = 42.1829218861219
$ws.Range("H8").Value = 17.58798119037683
$ws.Range("J8").Value = 9.009531971423042
$ws.Range("K8").Value = 9.303386565282217
$ws.Range("M8").Value = 17.1070211623634
$ws.Range("B9").Value = 10.40973720257323
$ws.Range("C9").Value = 4.610718671638072
$ws.Range("E9").Value = 18.18610851226791
$ws.Range("F9").Value = 44.32171678583015
$ws.Range("G9").Value = 42.38112712702839
$ws.Range("H9").Value = 17.51690861440956
$ws.Range("J9").Value = 8.965992080611297
$ws.Range("K9").Value = 9.63062537378029
$ws.Range("M9").Value = 17.31316387817158
$ws.Range("B10").Value = 10.7907568135259
$ws.Range("C10").Value = 4.846100009691321
$ws.Range("E10").Value = 18.34169440606017
$ws.Range("F10").Value = 44.61812507770134
$ws.Range("G10").Value = 42.58046694821107
$ws.Range("H10").Value = 17.47860941037231
$ws.Range("J10").Value = 8.936831716567728
$ws.Range("K10").Value = 9.874203384272842
$ws.Range("M10").Value = 17.47950848326939
$ws.Range("B11").Value = 10.96209085492015
$ws.Range("C11").Value = 4.948602724921225
$ws.Range("E11").Value = 18.41580433196267
$ws.Range("F11").Value = 44.7632192884261
$ws.Range("G11").Value = 42.68269649963
$ws.Range("H11").Value = 17.46422179998789
$ws.Range("J11").Value = 8.924174556429529
$ws.Range("K11").Value = 9.985196950115768
$ws.Range("M11").Value = 17.55819690544537
$ws.Range("B12").Value = 11.02661564374271
$ws.Range("C12").Value = 4.986744329566791
$ws.Range("E12").Value = 18.44432969387239
$ws.Range("F12").Value = 44.81960819963606
$ws.Range("G12").Value = 42.72305435486587
$ws.Range("H12").Value = 17.45921094195276
$ws.Range("J12").Value = 8.919468644948942
$ws.Range("K12").Value = 10.02721365216573
$ws.Range("M12").Value = 17.58840832653987
$ws.Range("B13").Value = 11.01273594737423
$ws.Range("C13").Value = 4.97856006239503
$ws.Range("E13").Value = 18.43816601034669
$ws.Range("F13").Value = 44.80740008043653
$ws.Range("G13").Value = 42.71428963977569
$ws.Range("H13").Value = 17.4602706506991
$ws.Range("J13").Value = 8.920478279868439
$ws.Range("K13").Value = 10.01816587804911
$ws.Range("M13").Value = 17.58188367839359
$ws.Range("B14").Value = 10.96740685635459
$ws.Range("C14").Value = 4.951754250020973
$ws.Range("E14").Value = 18.41814198506528
$ws.Range("F14").Value = 44.7678296609139
$ws.Range("G14").Value = 42.6859838765288
$ws.Range("H14").Value = 17.46380078179059
$ws.Range("J14").Value = 8.923785655057461
$ws.Range("K14").Value = 9.988654205492928
$ws.Range("M14").Value = 17.56067424101021
$ws.Range("B15").Value = 10.93959318453887
$ws.Range("C15").Value = 4.935246702581439
$ws.Range("E15").Value = 18.40593625778811
$ws.Range("F15").Value = 44.74377887104601
$ws.Range("G15").Value = 42.66885961477229
$ws.Range("H15").Value = 17.46602008178126
$ws.Range("J15").Value = 8.925822847833958
$ws.Range("K15").Value = 9.970574405093828
$ws.Range("M15").Value = 17.54773613273132
$ws.Range("B16").Value = 10.77951361420214
$ws.Range("C16").Value = 4.839307742214524
$ws.Range("E16").Value = 18.33691671861434
$ws.Range("F16").Value = 44.60884670955964
$ws.Range("G16").Value = 42.57401717930387
$ws.Range("H16").Value = 17.4796108324525
$ws.Range("J16").Value = 8.937671096604372
$ws.Range("K16").Value = 9.866949935984421
$ws.Range("M16").Value = 17.47442502755773
$ws.Range("B17").Value = 10.68074887057223
$ws.Range("C17").Value = 4.779268608838883
$ws.Range("E17").Value = 18.29541677157284
$ws.Range("F17").Value = 44.52867658992742
$ws.Range("G17").Value = 42.51878253984079
$ws.Range("H17").Value = 17.48872637255536
$ws.Range("J17").Value = 8.945095092697438
$ws.Range("K17").Value = 9.80339828918428
$ws.Range("M17").Value = 17.43020981143032
$ws.Range("B18").Value = 10.62375820160108
$ws.Range("C18").Value = 4.744306035484436
$ws.Range("E18").Value = 18.27186194287371
$ws.Range("F18").Value = 44.48353221095908
$ws.Range("G18").Value = 42.4881011338282
$ws.Range("H18").Value = 17.49425499444218
$ws.Range("J18").Value = 8.949422431042574
$ws.Range("K18").Value = 9.766865231676691
$ws.Range("M18").Value = 17.40506365155686
$ws.Range("B19").Value = 10.60443274652867
$ws.Range("C19").Value = 4.732395049645231
$ws.Range("E19").Value = 18.26394128564411
$ws.Range("F19").Value = 44.46841411801513
$ws.Range("G19").Value = 42.47790023213261
$ws.Range("H19").Value = 17.49617591141979
$ws.Range("J19").Value = 8.950897437427395
$ws.Range("K19").Value = 9.754500506329984
$ws.Range("M19").Value = 17.39659918847004
$ws.Range("B20").Value = 10.69128208258364
$ws.Range("C20").Value = 4.785704469638801
$ws.Range("E20").Value = 18.29980205289436
$ws.Range("F20").Value = 44.53711092811555
$ws.Range("G20").Value = 42.52454984753309
$ws.Range("H20").Value = 17.48772644114086
$ws.Range("J20").Value = 8.94429887283076
$ws.Range("K20").Value = 9.810161685735002
$ws.Range("M20").Value = 17.43488721627239
$ws.Range("B21").Value = 10.98073129069762
$ws.Range("C21").Value = 4.959646176458763
$ws.Range("E21").Value = 18.42401114582822
$ws.Range("F21").Value = 44.77941348968342
$ws.Range("G21").Value = 42.69425343254174
$ws.Range("H21").Value = 17.46275201782445
$ws.Range("J21").Value = 8.922811838263282
$ws.Range("K21").Value = 9.997323206491057
$ws.Range("M21").Value = 17.56689289920634
$ws.Range("B22").Value = 11.16780094902921
$ws.Range("C22").Value = 5.069392820555949
$ws.Range("E22").Value = 18.50787043342617
$ws.Range("F22").Value = 44.94617851331892
$ws.Range("G22").Value = 42.81474841908925
$ws.Range("H22").Value = 17.44897978722725
$ws.Range("J22").Value = 8.909276228115036
$ws.Range("K22").Value = 10.11954329710936
$ws.Range("M22").Value = 17.65556855340488
$ws.Range("B23").Value = 11.06817234964021
$ws.Range("C23").Value = 5.011183776727743
$ws.Range("E23").Value = 18.46287385624724
$ws.Range("F23").Value = 44.85641423343588
$ws.Range("G23").Value = 42.74956666542138
$ws.Range("H23").Value = 17.45609665559278
$ws.Range("J23").Value = 8.916454126222376
$ws.Range("K23").Value = 10.05433468762163
$ws.Range("M23").Value = 17.60802767955367
$ws.Range("B24").Value = 10.68652066467826
$ws.Range("C24").Value = 4.78279620269716
$ws.Range("E24").Value = 18.29781851880271
$ws.Range("F24").Value = 44.5332948132117
$ws.Range("G24").Value = 42.52193910190035
$ws.Range("H24").Value = 17.48817761311919
$ws.Range("J24").Value = 8.944658659443114
$ws.Range("K24").Value = 9.807103939668369
$ws.Range("M24").Value = 17.43277170728553
$ws.Range("B25").Value = 10.26813701856044
$ws.Range("C25").Value = 4.520103526871174
$ws.Range("E25").Value = 18.13188046831638
$ws.Range("F25").Value = 44.22183907318768
$ws.Range("G25").Value = 42.31804124858708
$ws.Range("H25").Value = 17.53369630639039
$ws.Range("J25").Value = 8.977272220096749
$ws.Range("K25").Value = 9.541371870843438
$ws.Range("M25").Value = 17.2547097091725

Write-Host "Updated 216 cells for 380 kV case"
